$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Country for each airline row (row 1 = header "country", rows 2-57 = data)
$countries = @("country", "Ireland", "Russia", "Argentina", "Mexico", "Canada", "France", "India", "New Zealand", "United States", "Italy", "Japan", "United States", "Austria", "Colombia", "United Kingdom", "Hong Kong", "China", "Germany", "Panama", "United States", "Egypt", "Israel", "Ethiopia", "Finland", "Indonesia", "Bahrain", "United States", "Spain", "Japan", "Kenya", "Netherlands", "South Korea", "Chile", "Germany", "Malaysia", "Pakistan", "Philippines", "Australia", "Morocco", "Sweden", "Saudi Arabia", "Singapore", "South Africa", "United States", "Sri Lanka", "Switzerland", "El Salvador", "Brazil", "Portugal", "Thailand", "Turkey", "United States", "United States", "Vietnam", "United States", "China")

for ($i = 0; $i -lt $countries.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $countries[$i]
}

# Resize column C to fit its new contents
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null

# Match the view state shown in the target workbook
$ws.Application.ActiveWindow.ScrollRow = 27
$ws.Range("C57").Select() | Out-Null
